$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bold the header row (A1:E1) ---
$ws.Range("A1:E1").Font.Bold = $true

# --- Remove existing hyperlinks on E2:E10 before rewriting their values ---
foreach ($hl in @($ws.Hyperlinks)) {
    $hl.Delete()
}

# --- Job listing data: Title, Company, Location, Salary, URL ---
$jobs = @(
    @('Junior Delphi Software Developer', 'Ernest Gordon Recruitment Limited', 'Hybrid remote in London', '£35,000 - £45,000 a year', 'https://uk.indeed.com/pagead/clk?mo=r&ad=-6NYlbfkN0BLTG9DKC8iUEjzj5GMXRvbrRPPLtQ49bnVbnj2Qcpi6Ewh8amSvyViVIVvvK1NOl7hiZPohhxO-TdX8jV39FkrpMiK9C8hZHtRCo_dMV8m-p2n-XFYiC3dllN0JP0XV78w2aKlygOu8Tbp6vtiYlP0FYxu2E06YjHNJeNGnXwBLtLfXFzKiyUG4On-paPT3R0fFGp4885R6YOHVL6VCoOPvZsgf0dkc_KGWpNUD1h1m_IlNkkO4XNBSHvq7RrZeeHd_3lqTdopA74eDOiNzkoobvtrXaCWaAzuyqr1eT6XEKBsHq4TYdDx2gEKUHYCEI3I7NMNOKqy_Gt721sCw_pe4KMXdftxWPr3Qt-YCOJPh1vzXjt7yhWqoEUU4b15efPajehfur_Dx-hD-O56mhLToh1I7XIk8_8FHn3C3sI9tJGbN5eK2X2VXTxnnRKpq0CX4b7hjIlWkLy14Z7FBqGWZ_H-x8Jn-YrjxEZ6_RNq9nzlw8WQmVSnf45BWEFFHCAvd8a92hia10u8HS_mwJoB5vjws84CvtfA0T9GEes3ig==&xkcb=SoAz-_M3bD5oxcRr0h0LbzkdCdPP&p=0&fvj=0&vjs=3'),
    @('Website Developer', 'Contra Agency', 'Temporarily Remote in London SE10', '£40,586 - £45,175 a year', 'https://uk.indeed.com/pagead/clk?mo=r&ad=-6NYlbfkN0A9KeT3j6INblOOMK0U8DFhLDoyAsNib7P09fbwmcp215U8nGCR2aXBvcYy7emKlCjwtYDPafTJzoM0GdZoRmsPES45zwLju91ayXln_Su9KaR9f3S1X_CtHRqESpkdguRX5gpzzvhiKxQ__jj7N3YUzyOD_e4stmreX4HAgxICna5Keqb9l6VWOagWobUxPe5sapNEXR1XJ-IW_URhfNjEji7EIgKADnc_bSRwrL4818fNbb0zJU9OSTNYUrWg0fMfM1ux0oAblqHMUAtDz8OM2jkDmreu2fsw2YsKK73J1Zg0cajcE3TqN2juS177L3KJ0UYjoNwj079gxL2cVWflfmbmUEl4-8W4QA5Ze0puIx_at4ourBsqu89KHHz-cJ7Gp8Fdy8rYYi9eW-RxZ4cQ5_hUw6Cw3DWh4xGW2IxVh9_Uu84ZBI4R9hYL5Jm7Zr_X2DUWlKnQCHKCh4QUdwKTIijXtTUJzbcUz3TVfgb_2SrtI1I4pnhw&xkcb=SoDC-_M3bD5oxcxr0h0KbzkdCdPP&p=1&fvj=1&vjs=3'),
    @('Delphi Software Developer (Remote)', 'Ernest Gordon Recruitment Limited', 'Remote in London', '£45,000 - £55,000 a year', 'https://uk.indeed.com/pagead/clk?mo=r&ad=-6NYlbfkN0BLTG9DKC8iUEjzj5GMXRvbrRPPLtQ49bnVbnj2Qcpi6Ewh8amSvyViVIVvvK1NOl5g0317UJAzIUtqgRtEvOQPB-t8BpmEIIPxQYPzV_yJckGzH2qz3vx6SLefP_uF1GrDOlOOLFBA5cdFc2haiyasTF4LforxBWUENCUX-XIR1Er1KBQV_xwxwJLt1GpKe4jgfx8Q8oHBguSVUO7xZ4xj3h6idzjnf8NuCykGlHPb0wzcjb96YwLnpqD_u2HtKJiFQ_YUceMtw6ZO8Gfr3JQv5CUktVknF8Yn8UlV9ajXLdGIeCvIaPYzMwiRXWm4cg5rqjKidGYJpVywn51m8ggE59cl_htWIE4zLtAnoEsQ6izEehH7RgF_rYZsoaromH-RkSm97KiW_XFDvZZxKq-gKoaq-jI9r6uHkIGQvJC4uL7xMnyxG2gedLcNbRMU5dQ7gtF7sGpfx253YqZ73YL_pLtYrlIVX-C5tRMsDPPNFwbsAZt1c1Yt06AGUEORDTbvp3Xa1I5i6Ee2TQj8RNeGxkrqDsoq_L6HZcpL_VzGqA==&xkcb=SoAv-_M3bD5oxfRr0h0LbzkdCdPP&p=3&fvj=0&vjs=3'),
    @('Software Developer', 'Ig Recruit Ltd', 'Hybrid remote in London', '£45,000 - £70,000 a year', 'https://uk.indeed.com/pagead/clk?mo=r&ad=-6NYlbfkN0DzFSXhtSQvn9oyRCOA6psFY_LWEYKJkh_kYNX4Ftteh7SDudwvmWFDGvVqCYY46KDuIXMvmGpvOr6r3eI6Efc9lqndbKp9fo_rOKkXQIhYoTVY4TvABX0Jy-aTsZ7_Zf65YHFmuCOe5wtM4FlyGIl-_mnDP-R2tcgrdgeXL0jh7e6QRiYc-IfyEtW78hINH3kkDACsI-GcVS-ZuXjtQjvHXChXLQkuHwgmC1lgYk47d0Xc6qUhDbWuBrYZTqbspRXf3sxD-f0DofM0UpM1qAwZcn8xWiHcWwRRuHz8CW8YXasYeCN3v0U8qcTwkEz-nZK7DdGyiICXMbIPuOI5cmlB_9Obc8G7fsWjz4hfmlkKssa0ycT4lUHEu5Btx8M5HLiRPBR9uphi3xxgJtSLKf_Pe_LZgFgFgJnWnUkP1Qai4afovvNtWDbkFVwbykwU_97N7yOWh9afus_o5ysIOa6MOKyxTj20Ocbqy98_0T8MjDVPxnaN__TT&xkcb=SoAG-_M3bD5oxfRr0h0JbzkdCdPP&p=4&fvj=1&vjs=3'),
    @('Delphi Software Developer (Hybrid)', 'Ernest Gordon Recruitment Limited', 'Hybrid remote in London', '£45,000 - £55,000 a year', 'https://uk.indeed.com/pagead/clk?mo=r&ad=-6NYlbfkN0BLTG9DKC8iUEjzj5GMXRvbrRPPLtQ49bnVbnj2Qcpi6Ewh8amSvyViVIVvvK1NOl5g0317UJAzIVCybnQxMbOUNjV32AFxvix0wEt4hkIY9mdTYZDYXmsdjvA9lGOMl8m1smYKzECwACTNwQC-EHSWKuQ0iuBeIV2GlMjeyFe686gIyEXSxzxXX1VdITPf94HfIZAhQ0SxPNysNIpakdvGKwzGuKvjoO4CJgkc0poggDUjWrBy3bUKkwMlr-bm_23ZBw25tN27iq2F2-GgkkW1ty2eaKtOrTIokt8db6pYrslKzQ_I9WcIBLKBlmoK2L1rjr00Tm19DA0hY9JGutQVGjADCTgnYD2i03RKFkihdrc0PjJb-JuKElWbw15mgKej-wbAmW_wB4pxBgpZGv1v43nxKTFcdRiC-ZzZr1IxAgACZQBpJqWYJjvFGiCk-50b5XQt1NlmLgQ1CLpaXkDnBSP0vdCO0sAoe6qNCNXx7XYxfSxmiUpB9ZsNXOZOzomEi-j2xHXNjcw3mg_OR2KZHWL970xjf48=&xkcb=SoBq-_M3bD5oxfxr0h0LbzkdCdPP&p=5&fvj=0&vjs=3'),
    @('Graduate Delphi Software Developer', 'Ernest Gordon Recruitment Limited', 'Hybrid remote in London', '£35,000 - £45,000 a year', 'https://uk.indeed.com/pagead/clk?mo=r&ad=-6NYlbfkN0BLTG9DKC8iUEjzj5GMXRvbrRPPLtQ49bnVbnj2Qcpi6Ewh8amSvyViVIVvvK1NOl7hiZPohhxO-S46WKZfQYYnhGZcpYpztB4Qipe9xy4-6DBJt5ytGu0Bxi0SSeCDkVO-cUXYb8e771fw1gDlm7K_xMyoPJwiCqhHYw7WW1zETPB_4ShCsL42eoiXCS2g95BMKyZdLrrNuEgJpEOTdQA45R7r7nzk07JFcUCrqUroXT3iZgC_GIGf4w1W2gXOk2EbZkbSpQHi_O4r-frdI6lbVcfr0mZqTjdlzy4BXRGJELw_hB31mZZ9WixvIWJC-0OfcDjoXjxToQ2KNiWC8WWgNFQ5PAZMi7aW387XK8NDflqywfxqNCHxBXfEfaCRUy4bqRJoAUdcWayjl8dnajVjuIgeJC7tLmBJSooYSaNyJi04E0iWyF9gh3XHkWlb8yLpv7FJF5elgXTT7v3nytywCJ36Jwcm6fwzeKO8Nb5WIEkSdCT93HWIKKzfV97KrZwO5rKr_8aESSjFXylijoLn98_kN6SZt3vdxHQCSd3x2KkhEVib0R__&xkcb=SoBD-_M3bD5oxfxr0h0JbzkdCdPP&p=7&fvj=0&vjs=3'),
    @('Software Developer - Your Return to Tech', 'Tech Returners', 'Temporarily Remote in London+2 locations', 'From £40,000 a year', 'https://uk.indeed.com/pagead/clk?mo=r&ad=-6NYlbfkN0AbMsrYsAi0JFHPLu_2mW_2EIGQHZf3spZOFbDQNHXQWmYvxYTDL77sEfykw7oHF0h1biF3zRwYbwETm6-5cxbE0wfjxa-ZNT-YEhj8B9himkV_JNb6VuMpkmPsSqn2T0XJ3qSM0bkMfXhtxNiuWXf22Kcc88kzNWA8CEbLSt_IFrN3opWHgfjP3FxODISygU3Cbvg-LNeg1fvGpELzwVli5k9PVV_QZcVMLohx8CXDEFsLS8DUuRZELksGTYQubs-1YOMSes5HXPDmM4mAQWRhUZ1MhUgjnUSclkmcwB8uloGNTTxAjtvv9UGeNeWPJGel9JNk6ToiQPZBs7hFOnI1A-jjWozehZI6mfXOLfxeNrTgek_L1PA8MGktdftgeOBxlKgqTPPYSdbW1DuEpyoj_y0AvnMZEwL0UlPbuIoj6pMciTiiMRgM_ea-f22GaTkneXhoz34eGTkcecjBLB0Skst6zpQwyf4v4-hNs75o5EqRj0s_mM4OWtLoE7xprkGIDFIdiRS4Hg==&xkcb=SoBQ-_M3bD5oxeRr0h0KbzkdCdPP&p=8&fvj=1&vjs=3'),
    @('Umbraco Developer - Fully Remote - £500+ per day - Outside I...', 'Edison Pope', 'Remote in London', 'From £500 a day', 'https://uk.indeed.com/pagead/clk?mo=r&ad=-6NYlbfkN0B62sP4HgT_bin6Q8ZuOedLyi0IXdvA6AE8k5EX4G5rP29lej8rB3Qex4moYYcYIbG1Z0Rg-2mE076IiCz8B8VfPWkv8uwOAvEBd1Q3qZvXQifJJh8okkaaA5M7y6XViUP1UEoGAQe3n20wEmrex6_YcCxbhfWoN39h7qfcwa6hjeCrS_AT-r1dAc2WE0IeHJcEylKfCft3IJeAi_pRh93pooTHsdizgIjhkfRmtCtw3trcSWpOjvyY4Ox26c00jViRzEc9KRiEZzLyOQZC-OX02hNUnLfM8Vr91k_7DmGjNTkhtZx0Y6jzjKJZyoO1EgtmGIFxTEBU4SiFeEjJYEttPgxp9iLZx0mAb4ZoZf80xo0E_Tl1O76EJKGRr5ULNoIwakpQiwLcyFnsfl4PvP_CPmcsRHLlgjAOtA5b4TfCtnAB3AciRQdNW_fzfmOdmRy1FEoR0C1RjmkYcAbaf8db1eZdnOUpEqcCT2A0m6YLftL1hmhQHyX1&xkcb=SoB5-_M3bD5oxeRr0h0IbzkdCdPP&p=9&fvj=1&vjs=3'),
    @('Junior JavaScript Developer', 'Human Exponent Inc.', 'London', '£30,000 - £40,000 a year', 'https://uk.indeed.com/company/Human-Exponent-Inc./jobs/Junior-Javascript-Developer-4ffb29857e1ae5d0?fccid=b4757f25c14c2b76&vjs=3'),
    @('Software Developer Internship', 'Inspiring Interns', 'London', '£18,000 - £19,000 a year', 'https://uk.indeed.com/pagead/clk?mo=r&ad=-6NYlbfkN0DzwreoEVFOkXRcFH1aYDKs9BsTnM4S4dXogYMNQuxRsWEvU3NmJBXmQipkwYVPqOyyWakXPnJku76CC3G11jOlJJl_Z-vD2e6CqEiCp4Xw5BmESzFcdir6NKaLvHdNE5m_HW0ROYcLmIm2F-BH_4hG_XB3HWa-WyT2T7s70oZyWPA5OdnFSdhBCOFQLiq8HX6O2PvWjintOy-bB3rwwWWHRitKYfS81uSFbsfQ1rv4XTjb4_YqfBwRyL8QHx0e00V9ydoge69SbNDyR9lI37fKetRQzFo080hexh9oJK14qcyfC_2ZFb3nZ7t6k13Tpk-W0we4BzLgH1S7bYlXmG0g28OtBS8YV0kDvJgghSSzCc17MURdMhfrvud_aBc2BViVMRY36VTmDeSIysca9pI_WWgqUiC-vk4yYkI7eh62x_NKLkdR4wkI3hg_qz9Ym-yzGM_MXAga4metIOJQq6OGmlse_Z4_eMb1mYGLh045mPHMHWgB2lsmIdpPeOxqVxU=&xkcb=SoCh-_M3bD5oxexr0h0LbzkdCdPP&p=11&fvj=1&vjs=3'),
    @('Software Development Internship', 'Softwire', 'London', '£28,000 a year', 'https://uk.indeed.com/rc/clk?jk=fc63bf03940b5708&fccid=ea94cafa08dade11&vjs=3'),
    @('Apprentice Software Developer', 'Softwire', 'London', '£24,000 a year', 'https://uk.indeed.com/rc/clk?jk=e97c14fa64f87964&fccid=ea94cafa08dade11&vjs=3')
)

$row = 2
foreach ($job in $jobs) {
    $ws.Cells.Item($row, 1).Value = $job[0]
    $ws.Cells.Item($row, 2).Value = $job[1]
    $ws.Cells.Item($row, 3).Value = $job[2]
    $ws.Cells.Item($row, 4).Value = $job[3]
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $job[4])
    $row = $row + 1
}

